$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94 - this shifts existing rows 94:106 down to 95:107
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45127
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100108
$ws.Range("H94").Value = "Tropicales y subtropicales"
$ws.Range("I94").Value = 100108007
$ws.Range("J94").Value = "Coco"
$ws.Range("K94").Value = "Sin especificar"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 25
$ws.Range("N94").Value = 36000
$ws.Range("O94").Value = 36000
$ws.Range("P94").Value = 36000
$ws.Range("Q94").Value = "$/malla 20 unidades"
$ws.Range("R94").Value = "Perú"
$ws.Range("S94").Value = 1800
$ws.Range("T94").Value = 20
